$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.732.52"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").Value = "2.548.38"
$ws.Range("E3").Value = "  +5.41%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'573.46"
$ws.Range("E5").Value = "  +2.53%  "
$ws.Range("D6").Value = "'148.40"
$ws.Range("E6").Value = "  +7.21%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.587"
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("D9").Value = "2.547.99"
$ws.Range("E9").Value = "  +5.52%  "
$ws.Range("D10").Value = "'0.106"
$ws.Range("E10").Value = "  +2.28%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("E12").Value = "  +1.96%  "
$ws.Range("D13").Value = "'0.357"
$ws.Range("E13").Value = "  +3.08%  "
$ws.Range("E14").Value = "  +9.17%  "
$ws.Range("D15").Value = "2.999.91"
$ws.Range("E15").Value = "  +5.50%  "
$ws.Range("D16").Value = "63.611.88"
$ws.Range("E16").Value = "  +2.55%  "
$ws.Range("E17").Value = "  +3.12%  "
$ws.Range("D18").Value = "2.544.17"
$ws.Range("E18").Value = "  +5.21%  "
$ws.Range("D19").Value = "'11.53"
$ws.Range("E19").Value = "  +4.01%  "
$ws.Range("D20").Value = "'341.71"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("E21").Value = "  +3.17%  "
$ws.Range("D22").Value = "'6.89"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'66.08"
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("E25").Value = "  -0.98%  "
$ws.Range("E26").Value = "  +3.60%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "'8.43"
$ws.Range("E28").Value = "  +1.36%  "
$ws.Range("D29").Value = "'1.43"
$ws.Range("E29").Value = "  +3.66%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").Value = "'7.04"
$ws.Range("E30").Value = "  +11.04%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0833"
$ws.Range("E31").Value = "  +6.83%  "
$ws.Range("E32").Value = "  +4.52%  "
$ws.Range("D33").Value = "'177.16"
$ws.Range("E33").Value = "  +3.30%  "
$ws.Range("D34").Value = "'1.61"
$ws.Range("E34").Value = "  +13.36%  "
$ws.Range("D35").Value = "'423.89"
$ws.Range("E35").Value = "  +11.76%  "
$ws.Range("E36").Value = "  +2.60%  "
$ws.Range("E37").Value = "  +3.04%  "
$ws.Range("D38").Value = "'4.45"
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("E40").Value = "  +4.81%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").Value = "'40.50"
$ws.Range("E42").Value = "  +3.72%  "
$ws.Range("D43").Value = "'155.97"
$ws.Range("E43").Value = "  +7.08%  "
$ws.Range("E44").Value = "  +3.74%  "
$ws.Range("D45").Value = "'20.91"
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("E46").Value = "  +4.25%  "
$ws.Range("D47").Value = "'0.0534"
$ws.Range("E47").Value = "  +3.18%  "
$ws.Range("D48").Value = "'0.0967"
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("E49").Value = "  +5.09%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'18.74"
$ws.Range("E50").Value = "  +4.65%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "'1.86"
$ws.Range("E51").Value = "  +11.00%  "
